$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'double[,]' 24,1
$arr[0,0] = 9.516480729190736
$arr[1,0] = 9.170275331694175
$arr[2,0] = 8.952472728208978
$arr[3,0] = 8.862542160691723
$arr[4,0] = 8.847542731773242
$arr[5,0] = 8.951264450970255
$arr[6,0] = 9.398275280443928
$arr[7,0] = 10.22784338344549
$arr[8,0] = 10.80220830364149
$arr[9,0] = 11.05475985586883
$arr[10,0] = 11.14906124350817
$arr[11,0] = 11.12881225169555
$arr[12,0] = 11.06254524801398
$arr[13,0] = 11.02177884345361
$arr[14,0] = 10.78552080502791
$arr[15,0] = 10.63828866698727
$arr[16,0] = 10.5527869720354
$arr[17,0] = 10.52369971070439
$arr[18,0] = 10.65404700959025
$arr[19,0] = 11.08204626844024
$arr[20,0] = 11.35396088979218
$arr[21,0] = 11.20957287717488
$arr[22,0] = 10.64692532296122
$arr[23,0] = 10.00917268723481
$ws.Range("B2:B25").Value = $arr

$arr = New-Object 'double[,]' 24,1
$arr[0,0] = 6.145569675762306
$arr[1,0] = 5.978901922792333
$arr[2,0] = 5.873289745482114
$arr[3,0] = 5.829467662113047
$arr[4,0] = 5.822144897733559
$arr[5,0] = 5.872701865985924
$arr[6,0] = 6.088805846638945
$arr[7,0] = 6.484986292711686
$arr[8,0] = 6.757320038897209
$arr[9,0] = 6.876796906513165
$arr[10,0] = 6.921380482161394
$arr[11,0] = 6.911808333974788
$arr[12,0] = 6.880478183748597
$arr[13,0] = 6.861200929279899
$arr[14,0] = 6.749420964820557
$arr[15,0] = 6.679699166830698
$arr[16,0] = 6.639183387809802
$arr[17,0] = 6.625395223770549
$arr[18,0] = 6.687164176477273
$arr[19,0] = 6.88969870471881
$arr[20,0] = 7.018211269704229
$arr[21,0] = 6.949982286849937
$arr[22,0] = 6.683790588622301
$arr[23,0] = 6.380973029697323
$ws.Range("C2:C25").Value = $arr

$arr = New-Object 'double[,]' 24,1
$arr[0,0] = 4.836593880964276
$arr[1,0] = 4.795374274760983
$arr[2,0] = 4.769532975394977
$arr[3,0] = 4.758873257366678
$arr[4,0] = 4.75709555491729
$arr[5,0] = 4.769389731554961
$arr[6,0] = 4.822494617429512
$arr[7,0] = 4.922209317477389
$arr[8,0] = 4.992511965249291
$arr[9,0] = 5.02379868257298
$arr[10,0] = 5.035541896473295
$arr[11,0] = 5.033017499027092
$arr[12,0] = 5.024766921545899
$arr[13,0] = 5.019699477398432
$arr[14,0] = 4.990452942181182
$arr[15,0] = 4.972330040193248
$arr[16,0] = 4.961841139402031
$arr[17,0] = 4.958278739483628
$arr[18,0] = 4.974266023574772
$arr[19,0] = 5.027193184834267
$arr[20,0] = 5.061173671832614
$arr[21,0] = 5.043094976288482
$arr[22,0] = 4.973390982536907
$arr[23,0] = 4.895734711070337
$ws.Range("D2:D25").Value = $arr

$arr = New-Object 'double[,]' 24,1
$arr[0,0] = 16.45662231652856
$arr[1,0] = 15.52626842187455
$arr[2,0] = 14.93045656831965
$arr[3,0] = 14.68175075304579
$arr[4,0] = 14.64010549543989
$arr[5,0] = 14.92712595775056
$arr[6,0] = 16.14108283617103
$arr[7,0] = 18.38985695548253
$arr[8,0] = 20.02410079230793
$arr[9,0] = 20.72552826646213
$arr[10,0] = 20.98513751507768
$arr[11,0] = 20.92949261587217
$arr[12,0] = 20.74700665777473
$arr[13,0] = 20.63444754702569
$arr[14,0] = 19.97741858707064
$arr[15,0] = 19.5636195268713
$arr[16,0] = 19.32165820945409
$arr[17,0] = 19.23905396751063
$arr[18,0] = 19.60807823426343
$arr[19,0] = 20.80076994803952
$arr[20,0] = 21.54526994555995
$arr[21,0] = 21.15110724077333
$arr[22,0] = 19.58799111030698
$arr[23,0] = 17.75155183517436
$ws.Range("E2:E25").Value = $arr

$arr = New-Object 'double[,]' 24,1
$arr[0,0] = 24.34937062023394
$arr[1,0] = 24.35783602769008
$arr[2,0] = 24.37059821565606
$arr[3,0] = 24.37769639679654
$arr[4,0] = 24.37898948570594
$arr[5,0] = 24.37068626856856
$arr[6,0] = 24.35071711371972
$arr[7,0] = 24.37172484708007
$arr[8,0] = 24.42395079512897
$arr[9,0] = 24.45569190829972
$arr[10,0] = 24.46885627261065
$arr[11,0] = 24.46597024946954
$arr[12,0] = 24.45675203145113
$arr[13,0] = 24.45125455560734
$arr[14,0] = 24.42203683650267
$arr[15,0] = 24.40615556839113
$arr[16,0] = 24.39777280098816
$arr[17,0] = 24.39506372327351
$arr[18,0] = 24.40776837536797
$arr[19,0] = 24.45942861275837
$arr[20,0] = 24.49986109966788
$arr[21,0] = 24.47767274343281
$arr[22,0] = 24.40703689632937
$arr[23,0] = 24.359582982698
$ws.Range("F2:F25").Value = $arr

$arr = New-Object 'double[,]' 24,1
$arr[0,0] = 7.344005520526261
$arr[1,0] = 7.344005520526261
$arr[2,0] = 7.344005520526261
$arr[3,0] = 7.344005520526261
$arr[4,0] = 7.344005520526261
$arr[5,0] = 7.344005520526261
$arr[6,0] = 7.344005520526261
$arr[7,0] = 7.344005520526261
$arr[8,0] = 7.344005520526261
$arr[9,0] = 7.344005520526261
$arr[10,0] = 7.344005520526261
$arr[11,0] = 7.344005520526261
$arr[12,0] = 7.344005520526261
$arr[13,0] = 7.344005520526261
$arr[14,0] = 7.344005520526261
$arr[15,0] = 7.344005520526261
$arr[16,0] = 7.344005520526261
$arr[17,0] = 7.344005520526261
$arr[18,0] = 7.344005520526261
$arr[19,0] = 7.344005520526261
$arr[20,0] = 7.344005520526261
$arr[21,0] = 7.344005520526261
$arr[22,0] = 7.344005520526261
$arr[23,0] = 7.344005520526261
$ws.Range("H2:H25").Value = $arr

$arr = New-Object 'double[,]' 24,1
$arr[0,0] = 8.667571513691255
$arr[1,0] = 8.417167355826319
$arr[2,0] = 8.260876708963771
$arr[3,0] = 8.196645470465585
$arr[4,0] = 8.185950188363023
$arr[5,0] = 8.260012523168339
$arr[6,0] = 8.581816195715779
$arr[7,0] = 9.18875105425378
$arr[8,0] = 9.643766359580635
$arr[9,0] = 9.886044951325937
$arr[10,0] = 9.976025520641786
$arr[11,0] = 9.956725541476844
$arr[12,0] = 9.893483242851399
$arr[13,0] = 9.854514801772407
$arr[14,0] = 9.627687477329481
$arr[15,0] = 9.492996367106738
$arr[16,0] = 9.429424487323439
$arr[17,0] = 9.407820200038252
$arr[18,0] = 9.504723621468671
$arr[19,0] = 9.912107166817567
$arr[20,0] = 10.17069553988727
$arr[21,0] = 10.03363346111814
$arr[22,0] = 9.499423301351545
$arr[23,0] = 9.027665718823961
$ws.Range("K2:K25").Value = $arr

$arr = New-Object 'double[,]' 24,1
$arr[0,0] = 21.75015368155668
$arr[1,0] = 21.81322534259342
$arr[2,0] = 21.85763764166023
$arr[3,0] = 21.87715984649082
$arr[4,0] = 21.88048730981546
$arr[5,0] = 21.85789516805771
$arr[6,0] = 21.77071714308888
$arr[7,0] = 21.64513309329788
$arr[8,0] = 21.58085853345163
$arr[9,0] = 21.55776100325154
$arr[10,0] = 21.54990236830025
$arr[11,0] = 21.55155530436786
$arr[12,0] = 21.55709664779564
$arr[13,0] = 21.56060664661999
$arr[14,0] = 21.5824920080862
$arr[15,0] = 21.59749443083587
$arr[16,0] = 21.6067012124893
$arr[17,0] = 21.6099175534514
$arr[18,0] = 21.59583756276172
$arr[19,0] = 21.55544488738954
$arr[20,0] = 21.53422322998391
$arr[21,0] = 21.54507442634601
$arr[22,0] = 21.59658482115925
$arr[23,0] = 21.67421395519014
$ws.Range("O2:O25").Value = $arr
